$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2360.4934
$ws.Range("I15").Value = 2360.4934
$ws.Range("K15").Value = 7081.4802
$ws.Range("M15").Value = -6912.4802

$ws.Range("H107").Value = 489.82352
$ws.Range("I107").Value = 624.9091
$ws.Range("J107").Value = 242.16667
$ws.Range("K107").Value = 624.9091
$ws.Range("L107").Value = 242.16667
$ws.Range("M107").Value = 1295.0909
$ws.Range("N107").Value = -4082.16667

$ws.Range("H137").Value = 4817.8965
$ws.Range("I137").Value = 1494.5294
$ws.Range("K137").Value = 4483.5882
$ws.Range("M137").Value = -1933.5882

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3421.896
$ws.Range("I32").Value = 2691.6028
$ws.Range("K32").Value = 2691.6028
$ws.Range("M32").Value = -2404.6028

$ws.Range("H57").Value = 10511
$ws.Range("I57").Value = 10511
$ws.Range("K57").Value = 10511
$ws.Range("M57").Value = -10027

$ws.Range("H61").Value = 4278.4443
$ws.Range("J61").Value = 6341.3
$ws.Range("L61").Value = 6341.3
$ws.Range("N61").Value = -6765.3

$ws.Range("H74").Value = 190345.03
$ws.Range("I74").Value = 293687.7
$ws.Range("K74").Value = 293687.7
$ws.Range("M74").Value = -292813.7

$ws.Range("H77").Value = 190345.03
$ws.Range("I77").Value = 293687.7
$ws.Range("K77").Value = 1468438.5
$ws.Range("M77").Value = -1464070.5

$ws.Range("H97").Value = 1646.174
$ws.Range("I97").Value = 1565.381
$ws.Range("J97").Value = 2494.5
$ws.Range("K97").Value = 1565.381
$ws.Range("L97").Value = 2494.5
$ws.Range("M97").Value = -1069.381
$ws.Range("N97").Value = -3486.5

$ws.Range("H110").Value = 4995.885
$ws.Range("I110").Value = 4343.7393
$ws.Range("K110").Value = 4343.7393
$ws.Range("M110").Value = -2298.7393

$ws.Range("H132").Value = 2015.4783
$ws.Range("I132").Value = 1332.8823
$ws.Range("K132").Value = 3998.6469
$ws.Range("M132").Value = -1468.6469

$ws.Range("H136").Value = 4278.4443
$ws.Range("J136").Value = 6341.3
$ws.Range("L136").Value = 19023.9
$ws.Range("N136").Value = -24123.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 62974.43
$ws.Range("I82").Value = 37245.5
$ws.Range("J82").Value = 73266
$ws.Range("K82").Value = 37245.5
$ws.Range("L82").Value = 73266
$ws.Range("M82").Value = -36862.5
$ws.Range("N82").Value = -74032

$ws.Range("H85").Value = 62974.43
$ws.Range("I85").Value = 37245.5
$ws.Range("J85").Value = 73266
$ws.Range("K85").Value = 37245.5
$ws.Range("L85").Value = 73266
$ws.Range("M85").Value = -35919.5
$ws.Range("N85").Value = -75918

$ws.Range("H86").Value = 3974.9375
$ws.Range("I86").Value = 3608.25
$ws.Range("J86").Value = 5075
$ws.Range("K86").Value = 3608.25
$ws.Range("L86").Value = 5075
$ws.Range("M86").Value = -2485.25
$ws.Range("N86").Value = -7321

$ws.Range("H89").Value = 3974.9375
$ws.Range("I89").Value = 3608.25
$ws.Range("J89").Value = 5075
$ws.Range("K89").Value = 18041.25
$ws.Range("L89").Value = 25375
$ws.Range("M89").Value = -12425.25
$ws.Range("N89").Value = -36607

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4473.364
$ws.Range("I31").Value = 3523.6316
$ws.Range("K31").Value = 3523.6316
$ws.Range("M31").Value = -3228.6316

$ws.Range("H34").Value = 4473.364
$ws.Range("I34").Value = 3523.6316
$ws.Range("K34").Value = 3523.6316
$ws.Range("M34").Value = -3321.6316

$ws.Range("H58").Value = 3652.0454
$ws.Range("I58").Value = 2735.182
$ws.Range("K58").Value = 2735.182
$ws.Range("M58").Value = -2532.182

$ws.Range("H132").Value = 15627881
$ws.Range("I132").Value = 18520710
$ws.Range("K132").Value = 55562130
$ws.Range("M132").Value = -55559600

$ws.Range("H134").Value = 6864.7646
$ws.Range("I134").Value = 6835.7856
$ws.Range("K134").Value = 20507.3568
$ws.Range("M134").Value = -17972.3568

$ws.Range("H136").Value = 3652.0454
$ws.Range("I136").Value = 2735.182
$ws.Range("K136").Value = 8205.545999999998
$ws.Range("M136").Value = -5655.545999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3502.25
$ws.Range("I3").Value = 2837
$ws.Range("K3").Value = 8511
$ws.Range("M3").Value = -8399

$ws.Range("H113").Value = 7304.778
$ws.Range("I113").Value = 495.66666
$ws.Range("J113").Value = 8666.6
$ws.Range("K113").Value = 1486.99998
$ws.Range("L113").Value = 25999.8
$ws.Range("M113").Value = 683.0000199999999
$ws.Range("N113").Value = -30339.8

$ws.Range("H133").Value = 5633
$ws.Range("I133").Value = 5633
$ws.Range("K133").Value = 16899
$ws.Range("M133").Value = -11839

$ws.Range("H134").Value = 1306.5714
$ws.Range("I134").Value = 1306.5714
$ws.Range("K134").Value = 3919.7142
$ws.Range("M134").Value = 1150.2858

$ws.Range("H136").Value = 995
$ws.Range("I136").Value = 995
$ws.Range("K136").Value = 2985
$ws.Range("M136").Value = 2115

$ws.Range("H139").Value = 3505.8572
$ws.Range("I139").Value = 2976.4
$ws.Range("K139").Value = 8929.200000000001
$ws.Range("M139").Value = -3789.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 727.8333
$ws.Range("I2").Value = 946.375
$ws.Range("K2").Value = 946.375
$ws.Range("M2").Value = -833.375

$ws.Range("H132").Value = 1991.037
$ws.Range("I132").Value = 1608.4736
$ws.Range("J132").Value = 2899.625
$ws.Range("K132").Value = 4825.4208
$ws.Range("L132").Value = 8698.875
$ws.Range("M132").Value = -2295.4208
$ws.Range("N132").Value = -13758.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 26875.773
$ws.Range("I40").Value = 29712.896
$ws.Range("J40").Value = 4746.2
$ws.Range("K40").Value = 29712.896
$ws.Range("L40").Value = 4746.2
$ws.Range("M40").Value = -29576.896
$ws.Range("N40").Value = -5018.2

$ws.Range("H61").Value = 1778.35
$ws.Range("I61").Value = 1620.4445
$ws.Range("K61").Value = 1620.4445
$ws.Range("M61").Value = -1418.4445

$ws.Range("H113").Value = 1778.35
$ws.Range("I113").Value = 1620.4445
$ws.Range("K113").Value = 1620.4445
$ws.Range("M113").Value = 549.5554999999999

$ws.Range("H136").Value = 5051.2607
$ws.Range("J136").Value = 5140.4287
$ws.Range("L136").Value = 15421.2861
$ws.Range("M136").Value = -20521.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 30000
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 30000
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H132").Value = 13337475
$ws.Range("I132").Value = 17548046
$ws.Range("K132").Value = 52644138
$ws.Range("M132").Value = -52641608

$ws.Range("H136").Value = 23258212
$ws.Range("I136").Value = 26316902
$ws.Range("K136").Value = 78950706
$ws.Range("M136").Value = -78948156
